$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2025/12/02 16:18"
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("G3").Value = "-"
